# Remove the hidden "Do you have TB?" slide (slide 30 in the deck order).
# This also causes PowerPoint to renumber/reflow the slide-number fields on
# the following slides automatically, matching the authoring diff, which
# shows slide30.xml removed and every subsequent slide shifted down by one
# position (slide31->30, slide32->31, ... slide36->35) together with the
# corresponding page-number text updates.
$p = $ppt.ActivePresentation
$p.Slides.Item(30).Delete()
